$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add SUM and AVERAGE formulas for rows 8-14 (months with no data yet),
# mirroring the pattern already present in rows 3-7.
for ($r = 8; $r -le 14; $r++) {
    $ws.Range("G$r").Formula = "=SUM(C${r}:E${r})"
    $ws.Range("H$r").Formula = "=AVERAGE(C${r}:E${r})"
}

# Update the active cell selection on the sheet view from M15 to H15.
$ws.Range("H15").Select()
